# "Generate Report for Handback"
#
# The localization-status report is refreshed after a handback: the
# "Ready for handoff" status becomes "Handed back: in sync with en-US" on
# the Overview sheet and on each per-locale sheet's Status column, and the
# per-locale sheets gain the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" values (with a hyperlink on the target-file
# cell) for each row.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/254d1d5c5964a7ad5183e2a4fa816083943b9e03/e2e/a.md"

# ---------------------------------------------------------------------
# Overview sheet: refresh the per-locale status cells for both rows.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack
$overview.Range("E3").Value = $statusHandedBack
$overview.Range("F3").Value = $statusHandedBack

$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# Per-locale sheets ("zh-cn", "de-de"): Status column + handback columns.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$zhHandbackTime = "2016-08-17 20:35:20"
$deHandbackTime = "2016-08-17 20:35:27"

foreach ($row in 2, 3) {
    # Status column (C) now reflects the handback.
    $zhcn.Range("C$row").Value = $statusHandedBack
    $dede.Range("C$row").Value = $statusHandedBack

    # "Latest Target File" (I) -- hyperlinked to the source markdown file.
    $zhcn.Hyperlinks.Add($zhcn.Range("I$row"), $aMdUrl, "", "", "a.md")
    $zhcn.Range("I$row").Font.Underline = $true
    $zhcn.Range("I$row").Font.Color = 15570276

    $dede.Hyperlinks.Add($dede.Range("I$row"), $aMdUrl, "", "", "a.md")
    $dede.Range("I$row").Font.Underline = $true
    $dede.Range("I$row").Font.Color = 15570276

    # "Latest Handback File" (J)
    $zhcn.Range("J$row").Value = $zhHandbackFile
    $dede.Range("J$row").Value = $deHandbackFile

    # "Latest Handback DateTime" (K)
    $zhcn.Range("K$row").Value = $zhHandbackTime
    $dede.Range("K$row").Value = $deHandbackTime
}

$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(10).ColumnWidth = 29.1

$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(10).ColumnWidth = 29.1
